$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.593.99"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "3.475.88"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.24"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.26"
$ws.Range("E6").Value = "  -1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  +2.63%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "3.476.32"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("E11").Value = "  -2.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -2.99%  "

$ws.Range("D13").Value = "4.087.53"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.20"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.132"
$ws.Range("E15").Value = "  -2.69%  "

$ws.Range("D16").Value = "67.616.85"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").Value = "  -2.68%  "

$ws.Range("D18").Value = "3.468.67"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -3.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  -3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.23"
$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.89"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.20"
$ws.Range("E25").Value = "  -2.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.535"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  -4.67%  "

$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  -4.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.50"
$ws.Range("E32").Value = "  +2.02%  "

$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.38"
$ws.Range("E34").Value = "  -4.99%  "

$ws.Range("E35").Value = "  -3.34%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("E37").Value = "  -5.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.46"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.888"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.78"
$ws.Range("E40").Value = "  +4.65%  "

$ws.Range("E41").Value = "  -4.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  -4.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.61"
$ws.Range("E43").Value = "  -5.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.54"
$ws.Range("E44").Value = "  -4.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0709"
$ws.Range("E45").Value = "  -3.71%  "

$ws.Range("D46").Value = "2.721.73"
$ws.Range("E46").Value = "  -5.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.84"
$ws.Range("E47").Value = "  -3.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.53"
$ws.Range("E48").Value = "  -2.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0297"
$ws.Range("E49").Value = "  -2.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "328.14"
$ws.Range("E50").Value = "  -6.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.04"
$ws.Range("E51").Value = "  -3.49%  "
